# ---------------------------------------------------------------------------
# Move develop template.xlsx and template_type.xlsx to 1.2.0.1
# Appends 40 new template-type rows (1704-1743) for: IDP auth-factor labels,
# VID-card-download notification templates, and supervisor-reject templates
# across the languages already present in the sheet (eng/fra/ara/hin/kan/tam).
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D on every new row repeats the existing "TRUE" (is_active) text cell,
# copied (value + style) from D1703 so it reuses the same shared string / style.
$ws.Cells.Item(1703, 4).Copy()
for ($r = 1704; $r -le 1743; $r++) {
    $ws.Cells.Item($r, 4).PasteSpecial(-4104)
}

# Rows 1704-1706 (the mosip.idp.* property keys) are highlighted with the
# JetBrains-Mono "changed" style used elsewhere in the sheet (same font/size,
# but the green-ish FF6A8759 color instead of FF067D17).
$ws.Cells.Item(1699, 3).Copy()
$ws.Range("B1704:B1706").PasteSpecial(-4122)
$ws.Range("B1704:B1706").Font.Color = 5867370

$ws.Cells.Item(1704, 1).Value = "eng"
$ws.Cells.Item(1704, 2).Value = "mosip.idp.otp.template.property"
$ws.Cells.Item(1704, 3).Value = "OTP"
$ws.Cells.Item(1705, 1).Value = "eng"
$ws.Cells.Item(1705, 2).Value = "mosip.idp.biometrics.template.property"
$ws.Cells.Item(1705, 3).Value = "Biometrics"
$ws.Cells.Item(1706, 1).Value = "eng"
$ws.Cells.Item(1706, 2).Value = "mosip.idp.unknown.authentication.template.property"
$ws.Cells.Item(1706, 3).Value = "Unknown"

$ws.Cells.Item(1707, 1).Value = "eng"
$ws.Cells.Item(1708, 1).Value = "eng"
$ws.Cells.Item(1709, 1).Value = "eng"
$ws.Cells.Item(1710, 1).Value = "eng"
$ws.Cells.Item(1711, 1).Value = "eng"
$ws.Cells.Item(1712, 1).Value = "eng"
$ws.Cells.Item(1713, 1).Value = "eng"
$ws.Cells.Item(1714, 1).Value = "eng"
$ws.Cells.Item(1715, 1).Value = "eng"
$ws.Cells.Item(1707, 2).Value = "mosip.full.name.template.property"
$ws.Cells.Item(1708, 2).Value = "mosip.date.of.birth.template.property"
$ws.Cells.Item(1709, 2).Value = "mosip.uin.template.property"
$ws.Cells.Item(1710, 2).Value = "mosip.perpetual.vid.template.property"
$ws.Cells.Item(1711, 2).Value = "mosip.phone.template.property"
$ws.Cells.Item(1712, 2).Value = "mosip.email.template.property"
$ws.Cells.Item(1713, 2).Value = "mosip.address.template.property"
$ws.Cells.Item(1714, 2).Value = "mosip.gender.template.property"
$ws.Cells.Item(1715, 2).Value = "mosip.defualt.template.property"
$ws.Cells.Item(1707, 3).Value = "Full Name"
$ws.Cells.Item(1708, 3).Value = "Date Of Birth"
$ws.Cells.Item(1709, 3).Value = "UIN"
$ws.Cells.Item(1710, 3).Value = "Perpetual VID"
$ws.Cells.Item(1711, 3).Value = "Phone"
$ws.Cells.Item(1713, 3).Value = "Address"
$ws.Cells.Item(1714, 3).Value = "Gender"
$ws.Cells.Item(1715, 3).Value = "Defualt"
$ws.Cells.Item(1712, 3).Value = "Email"

$ws.Cells.Item(1716, 1).Value = "eng"
$ws.Cells.Item(1716, 2).Value = "vid-card-type"
$ws.Cells.Item(1716, 3).Value = "Vid Card Type"
$ws.Cells.Item(1717, 1).Value = "eng"
$ws.Cells.Item(1717, 2).Value = "vid-card-download-request-received-email-subject"
$ws.Cells.Item(1717, 3).Value = "Request received email subject to download my VID card"
$ws.Cells.Item(1718, 1).Value = "eng"
$ws.Cells.Item(1718, 2).Value = "vid-card-download-success-email-subject"
$ws.Cells.Item(1718, 3).Value = "Success email subject to download my VID card"
$ws.Cells.Item(1719, 1).Value = "eng"
$ws.Cells.Item(1719, 2).Value = "vid-card-download-failure-email-subject"
$ws.Cells.Item(1719, 3).Value = "Failure email subject to download my VID card"
$ws.Cells.Item(1720, 1).Value = "eng"
$ws.Cells.Item(1720, 2).Value = "vid-card-download-request-received-email-content"
$ws.Cells.Item(1720, 3).Value = "Request received email to download my VID card"
$ws.Cells.Item(1721, 1).Value = "eng"
$ws.Cells.Item(1721, 2).Value = "vid-card-download-success-email-content"
$ws.Cells.Item(1721, 3).Value = "Success email to download my VID card"
$ws.Cells.Item(1722, 1).Value = "eng"
$ws.Cells.Item(1722, 2).Value = "vid-card-download-failure-email-content"
$ws.Cells.Item(1722, 3).Value = "Failure email to download my VID card"
$ws.Cells.Item(1723, 1).Value = "eng"
$ws.Cells.Item(1723, 2).Value = "vid-card-download-request-received_SMS"
$ws.Cells.Item(1723, 3).Value = "Request received sms to download my VID card"
$ws.Cells.Item(1724, 1).Value = "eng"
$ws.Cells.Item(1724, 2).Value = "vid-card-download-success_SMS"
$ws.Cells.Item(1724, 3).Value = "Success sms to download my VID card"
$ws.Cells.Item(1725, 1).Value = "eng"
$ws.Cells.Item(1725, 2).Value = "vid-card-download-failure_SMS"
$ws.Cells.Item(1725, 3).Value = "Failure sms to download my VID card"
$ws.Cells.Item(1726, 1).Value = "eng"
$ws.Cells.Item(1726, 2).Value = "RPR_SUP_REJECT_EMAIL"
$ws.Cells.Item(1726, 3).Value = "Template for Supervisor Reject Email"
$ws.Cells.Item(1727, 1).Value = "eng"
$ws.Cells.Item(1727, 2).Value = "RPR_SUP_REJECT_SMS"
$ws.Cells.Item(1727, 3).Value = "Template for Supervisor Reject SMS"
$ws.Cells.Item(1728, 1).Value = "fra"
$ws.Cells.Item(1728, 2).Value = "RPR_SUP_REJECT_EMAIL"
$ws.Cells.Item(1728, 3).Value = "Template for Supervisor Reject Email"
$ws.Cells.Item(1729, 1).Value = "fra"
$ws.Cells.Item(1729, 2).Value = "RPR_SUP_REJECT_SMS"
$ws.Cells.Item(1729, 3).Value = "Template for Supervisor Reject SMS"
$ws.Cells.Item(1730, 1).Value = "ara"
$ws.Cells.Item(1730, 2).Value = "RPR_SUP_REJECT_EMAIL"
$ws.Cells.Item(1730, 3).Value = "Template for Supervisor Reject Email"
$ws.Cells.Item(1731, 1).Value = "ara"
$ws.Cells.Item(1731, 2).Value = "RPR_SUP_REJECT_SMS"
$ws.Cells.Item(1731, 3).Value = "Template for Supervisor Reject SMS"
$ws.Cells.Item(1732, 1).Value = "hin"
$ws.Cells.Item(1732, 2).Value = "RPR_SUP_REJECT_EMAIL"
$ws.Cells.Item(1732, 3).Value = "Template for Supervisor Reject Email"
$ws.Cells.Item(1733, 1).Value = "hin"
$ws.Cells.Item(1733, 2).Value = "RPR_SUP_REJECT_SMS"
$ws.Cells.Item(1733, 3).Value = "Template for Supervisor Reject SMS"
$ws.Cells.Item(1734, 1).Value = "kan"
$ws.Cells.Item(1734, 2).Value = "RPR_SUP_REJECT_EMAIL"
$ws.Cells.Item(1734, 3).Value = "Template for Supervisor Reject Email"
$ws.Cells.Item(1735, 1).Value = "kan"
$ws.Cells.Item(1735, 2).Value = "RPR_SUP_REJECT_SMS"
$ws.Cells.Item(1735, 3).Value = "Template for Supervisor Reject SMS"
$ws.Cells.Item(1736, 1).Value = "tam"
$ws.Cells.Item(1736, 2).Value = "RPR_SUP_REJECT_EMAIL"
$ws.Cells.Item(1736, 3).Value = "Template for Supervisor Reject Email"
$ws.Cells.Item(1737, 1).Value = "tam"
$ws.Cells.Item(1737, 2).Value = "RPR_SUP_REJECT_SMS"
$ws.Cells.Item(1737, 3).Value = "Template for Supervisor Reject SMS"
$ws.Cells.Item(1738, 1).Value = "eng"
$ws.Cells.Item(1738, 2).Value = "RPR_SUP_REJECT_EMAIL_SUBJECT"
$ws.Cells.Item(1738, 3).Value = "Template for Supervisor Reject Email Subject"
$ws.Cells.Item(1739, 1).Value = "fra"
$ws.Cells.Item(1739, 2).Value = "RPR_SUP_REJECT_EMAIL_SUBJECT"
$ws.Cells.Item(1739, 3).Value = "Template for Supervisor Reject Email Subject"
$ws.Cells.Item(1740, 1).Value = "ara"
$ws.Cells.Item(1740, 2).Value = "RPR_SUP_REJECT_EMAIL_SUBJECT"
$ws.Cells.Item(1740, 3).Value = "Template for Supervisor Reject Email Subject"
$ws.Cells.Item(1741, 1).Value = "hin"
$ws.Cells.Item(1741, 2).Value = "RPR_SUP_REJECT_EMAIL_SUBJECT"
$ws.Cells.Item(1741, 3).Value = "Template for Supervisor Reject Email Subject"
$ws.Cells.Item(1742, 1).Value = "kan"
$ws.Cells.Item(1742, 2).Value = "RPR_SUP_REJECT_EMAIL_SUBJECT"
$ws.Cells.Item(1742, 3).Value = "Template for Supervisor Reject Email Subject"
$ws.Cells.Item(1743, 1).Value = "tam"
$ws.Cells.Item(1743, 2).Value = "RPR_SUP_REJECT_EMAIL_SUBJECT"
$ws.Cells.Item(1743, 3).Value = "Template for Supervisor Reject Email Subject"

# Leave the cursor where the author left it after typing the last row.
$ws.Range("F1730").Select()
